# "pruebas hasta conbenio de recaudos full"
#
# 1. Fix typo'd IP literal in the credentials block (row 14, col B):
#    "1092.168.0.1" -> "192.168.0.1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "192.168.0.1"

# 2. Append a new row (19) listing the client names used by the "recaudos"
#    convenios: edeq / tigoUne / movistar, under a new header label.
$ws.Range("A19").Value = "Nombre lista cliente"
$ws.Range("B19").Value = "edeq"
$ws.Range("C19").Value = "tigoUne"
$ws.Range("D19").Value = "movistar"

# 3. Leave the sheet scrolled/selected where the author left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("F15").Select() | Out-Null
